$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The alcohol measurement data had a redundant/duplicate "Total" column; the
# last data column (M, 13th) is removed so the following column shifts left
# to become the new column M.
$ws.Columns.Item(13).Delete() | Out-Null

# Excel leaves the selection where the deleted column used to be.
$ws.Range("M1").Select() | Out-Null
